$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("F7").Value = 'Ol. Grudziadz'
$ws.Range("G7").Value = [double]3
$ws.Range("H7").Value = 'Sandecja Nowy S.'
$ws.Range("I7").Value = [double]1
$ws.Range("J7").Value = [double]2.49
$ws.Range("K7").Value = '22/07/2023 14:42'
$ws.Range("L7").Value = [double]2.39
$ws.Range("M7").Value = '22/07/2023 17:36'
$ws.Range("N7").Value = [double]3.26
$ws.Range("O7").Value = '22/07/2023 14:42'
$ws.Range("P7").Value = [double]3.35
$ws.Range("Q7").Value = '22/07/2023 17:36'
$ws.Range("R7").Value = [double]2.64
$ws.Range("S7").Value = '22/07/2023 14:42'
$ws.Range("T7").Value = [double]2.82
$ws.Range("U7").Value = '22/07/2023 17:36'
$ws.Range("V7").Value = 'https://www.betexplorer.com/football/poland/division-2/ol-grudziadz-sandecja-nowy-s/GUiPyHqO/'

# Row 8
$ws.Range("F8").Value = 'Wisla Pulawy'
$ws.Range("G8").Value = [double]1
$ws.Range("H8").Value = 'Chojniczanka'
$ws.Range("I8").Value = [double]0
$ws.Range("J8").Value = [double]2.18
$ws.Range("K8").Value = '22/07/2023 14:42'
$ws.Range("L8").Value = [double]2.46
$ws.Range("M8").Value = '22/07/2023 17:06'
$ws.Range("N8").Value = [double]3.31
$ws.Range("O8").Value = '22/07/2023 14:42'
$ws.Range("P8").Value = [double]3.45
$ws.Range("Q8").Value = '22/07/2023 17:35'
$ws.Range("R8").Value = [double]3.05
$ws.Range("S8").Value = '22/07/2023 14:42'
$ws.Range("T8").Value = [double]2.6
$ws.Range("U8").Value = '22/07/2023 17:06'
$ws.Range("V8").Value = 'https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-chojniczanka/rc7er10t/'

# Row 134
$ws.Range("F134").Value = 'Wisla Pulawy'
$ws.Range("G134").Value = [double]2
$ws.Range("H134").Value = 'GKS Jastrzebie'
$ws.Range("I134").Value = [double]2
$ws.Range("J134").Value = [double]1.83
$ws.Range("K134").Value = '02/11/2023 06:12'
$ws.Range("L134").Value = [double]2.16
$ws.Range("M134").Value = '03/11/2023 17:25'
$ws.Range("N134").Value = [double]3.48
$ws.Range("O134").Value = '02/11/2023 06:12'
$ws.Range("P134").Value = [double]3.51
$ws.Range("Q134").Value = '03/11/2023 17:26'
$ws.Range("R134").Value = [double]3.7
$ws.Range("S134").Value = '02/11/2023 06:12'
$ws.Range("T134").Value = [double]3.1
$ws.Range("U134").Value = '03/11/2023 17:25'
$ws.Range("V134").Value = 'https://www.betexplorer.com/football/poland/division-2/wisla-pulawy-gks-jastrzebie/QyyJBxqo/'

# Row 135
$ws.Range("F135").Value = 'Polonia Bytom'
$ws.Range("G135").Value = [double]2
$ws.Range("H135").Value = 'Olimpia Elblag'
$ws.Range("I135").Value = [double]1
$ws.Range("J135").Value = [double]2.48
$ws.Range("K135").Value = '02/11/2023 06:12'
$ws.Range("L135").Value = [double]2.12
$ws.Range("M135").Value = '03/11/2023 17:24'
$ws.Range("N135").Value = [double]2.99
$ws.Range("O135").Value = '02/11/2023 06:12'
$ws.Range("P135").Value = [double]3.17
$ws.Range("Q135").Value = '03/11/2023 16:12'
$ws.Range("R135").Value = [double]2.65
$ws.Range("S135").Value = '02/11/2023 06:12'
$ws.Range("T135").Value = [double]3.52
$ws.Range("U135").Value = '03/11/2023 17:24'
$ws.Range("V135").Value = 'https://www.betexplorer.com/football/poland/division-2/polonia-bytom-olimpia-elblag/z5wjEV4q/'

# Row 148
$ws.Range("F148").Value = 'Lech Poznan II'
$ws.Range("G148").Value = [double]1
$ws.Range("H148").Value = 'Polonia Bytom'
$ws.Range("I148").Value = [double]0
$ws.Range("J148").Value = [double]2.89
$ws.Range("K148").Value = '11/11/2023 01:13'
$ws.Range("L148").Value = [double]2.98
$ws.Range("M148").Value = '12/11/2023 12:51'
$ws.Range("N148").Value = [double]3.36
$ws.Range("O148").Value = '11/11/2023 01:13'
$ws.Range("P148").Value = [double]3.55
$ws.Range("Q148").Value = '12/11/2023 12:51'
$ws.Range("R148").Value = [double]2.17
$ws.Range("S148").Value = '11/11/2023 01:13'
$ws.Range("T148").Value = [double]2.21
$ws.Range("U148").Value = '12/11/2023 12:51'
$ws.Range("V148").Value = 'https://www.betexplorer.com/football/poland/division-2/lech-poznan-polonia-bytom/8IybC9Ze/'

# Row 149
$ws.Range("F149").Value = 'Stezyca'
$ws.Range("G149").Value = [double]1
$ws.Range("H149").Value = 'Stomil Olsztyn'
$ws.Range("I149").Value = [double]0
$ws.Range("J149").Value = [double]1.95
$ws.Range("K149").Value = '11/11/2023 01:13'
$ws.Range("L149").Value = [double]1.9
$ws.Range("M149").Value = '12/11/2023 12:39'
$ws.Range("N149").Value = [double]3.28
$ws.Range("O149").Value = '11/11/2023 01:13'
$ws.Range("P149").Value = [double]3.33
$ws.Range("Q149").Value = '12/11/2023 12:39'
$ws.Range("R149").Value = [double]3.35
$ws.Range("S149").Value = '11/11/2023 01:13'
$ws.Range("T149").Value = [double]4.08
$ws.Range("U149").Value = '12/11/2023 12:39'
$ws.Range("V149").Value = 'https://www.betexplorer.com/football/poland/division-2/stezyca-stomil-olsztyn/dhXiYjdF/'

# Row 150
$ws.Range("F150").Value = 'Zaglebie II'
$ws.Range("G150").Value = [double]4
$ws.Range("H150").Value = 'S. Wola'
$ws.Range("I150").Value = [double]0
$ws.Range("J150").Value = [double]2.3
$ws.Range("K150").Value = '11/11/2023 01:13'
$ws.Range("L150").Value = [double]2.72
$ws.Range("M150").Value = '12/11/2023 12:51'
$ws.Range("N150").Value = [double]3.19
$ws.Range("O150").Value = '11/11/2023 01:13'
$ws.Range("P150").Value = [double]3.15
$ws.Range("Q150").Value = '12/11/2023 12:51'
$ws.Range("R150").Value = [double]2.81
$ws.Range("S150").Value = '11/11/2023 01:13'
$ws.Range("T150").Value = [double]2.59
$ws.Range("U150").Value = '12/11/2023 12:51'
$ws.Range("V150").Value = 'https://www.betexplorer.com/football/poland/division-2/zaglebie-stal-stalowa-wola/zqWeXABL/'

# Row 152 (new row) - copy styles from row 151 for A and E columns first
$ws.Range("A151").Copy() | Out-Null
$ws.Range("A152").PasteSpecial(-4122) | Out-Null
$ws.Range("E151").Copy() | Out-Null
$ws.Range("E152").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A152").Value = [double]151
$ws.Range("B152").Value = 'poland'
$ws.Range("C152").Value = 'division-2'
$ws.Range("D152").Value = '2023-2024'
$ws.Range("E152").Value = [double]45243.76041666666
$ws.Range("F152").Value = 'Skra'
$ws.Range("G152").Value = [double]4
$ws.Range("H152").Value = 'KKS Kalisz'
$ws.Range("I152").Value = [double]1
$ws.Range("J152").Value = [double]2.5
$ws.Range("K152").Value = '12/11/2023 06:42'
$ws.Range("L152").Value = [double]2.25
$ws.Range("M152").Value = '13/11/2023 18:14'
$ws.Range("N152").Value = [double]3
$ws.Range("O152").Value = '12/11/2023 06:42'
$ws.Range("P152").Value = [double]3.2
$ws.Range("Q152").Value = '13/11/2023 18:14'
$ws.Range("R152").Value = [double]2.63
$ws.Range("S152").Value = '12/11/2023 06:42'
$ws.Range("T152").Value = [double]3.19
$ws.Range("U152").Value = '13/11/2023 18:14'
$ws.Range("V152").Value = 'https://www.betexplorer.com/football/poland/division-2/skra-czestochowa-kks-kalisz/2HIlZWt9/'
